{"js": "// The edit: the author placed the cursor at the end of the \"Siiiii\"\n// paragraph (just before the hidden _GoBack bookmark), pressed Enter and\n// typed \"As\u00ed es\" -- producing a new paragraph \"As\u00ed es\" right after\n// \"Siiiii\", with the _GoBack bookmark carried forward onto the (now)\n// trailing empty paragraph, and Word's proofer wrapping the now-isolated\n// \"Siiiii\" word with spellStart/spellEnd proofErr marks (just like the\n// pre-existing misspelled words \"Git\" and \"Agrregar\" in this document).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// \"Siiiii\" is the 3rd paragraph (index 2).\nconst siiiiiPara = paragraphs.items[2];\n\n// 1) Split off a new paragraph \"As\u00ed es\" right after \"Siiiii\".\nsiiiiiPara.insertParagraph(\"As\u00ed es\", Word.InsertLocation.after);\nawait context.sync();\n\n// 2) The _GoBack bookmark used to sit in the \"Siiiii\" paragraph; move it\n//    onto the trailing (now last) paragraph, mirroring where it ends up\n//    after the in-place edit/save.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphsAfterSplit = body.paragraphs;\nparagraphsAfterSplit.load(\"items\");\nawait context.sync();\nconst trailingPara = paragraphsAfterSplit.items[paragraphsAfterSplit.items.length - 1];\ntrailingPara.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Wrap the now-isolated \"Siiiii\" run with spellcheck proofErr marks\n//    (spellStart/spellEnd), matching Word's automatic proofing markup\n//    already present around the other misspelled words in this doc.\nconst siiiiiOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Siiiii</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nsiiiiiPara.getRange().insertOoxml(siiiiiOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The edit: the author placed the cursor at the end of the \"Siiiii\"\n# paragraph (just before the hidden _GoBack bookmark), pressed Enter and\n# typed \"As\u00ed es\" -- producing a new paragraph \"As\u00ed es\" right after\n# \"Siiiii\", with the _GoBack bookmark carried forward onto the (now)\n# trailing empty paragraph, and Word's proofer wrapping the now-isolated\n# \"Siiiii\" word with spellStart/spellEnd proofErr marks (just like the\n# pre-existing misspelled words \"Git\" and \"Agrregar\" in this document).\n\n$d = $word.ActiveDocument\n\n# \"Siiiii\" is the 3rd paragraph.\n$siiiiiPara = $d.Paragraphs.Item(3)\n\n# 1) Split off a new paragraph \"As\u00ed es\" right after \"Siiiii\" by inserting\n#    a paragraph break + text just before the paragraph mark (i.e. right\n#    where the hidden _GoBack bookmark sits).\n$splitPos = $siiiiiPara.Range.End - 1\n$splitRange = $d.Range($splitPos, $splitPos)\n$splitRange.InsertAfter([char]13 + \"As\u00ed es\")\n\n# 2) The _GoBack bookmark used to sit in the \"Siiiii\" paragraph; move it\n#    onto the trailing (now last) paragraph, mirroring where it ends up\n#    after the in-place edit/save.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$trailingPara = $d.Paragraphs.Last\n$d.Bookmarks.Add(\"_GoBack\", $trailingPara.Range)\n\n# 3) Wrap the now-isolated \"Siiiii\" run with spellcheck proofErr marks\n#    (spellStart/spellEnd), matching Word's automatic proofing markup\n#    already present around the other misspelled words in this doc.\n$siiiiiXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Siiiii</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n$siiiiiPara.Range.InsertXML($siiiiiXml)\n"}
